$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date string (slash -> dash)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$r]
}

# Map of row -> @{ D = ..; E = ..; H = ..; G = .. } for cells whose numeric values change
$values = @{
    3  = @{ D = 1; G = 1 }
    4  = @{ D = 1; E = 1; H = 0 }
    5  = @{ D = 1; E = 1; H = 0 }
    6  = @{ D = 1; E = 1; H = 0 }
    7  = @{ D = 1; E = 1; H = 0 }
    10 = @{ D = 1; E = 1; H = 0 }
    12 = @{ D = 1; E = 1; H = 0 }
    14 = @{ D = 1; E = 1; H = 0 }
    20 = @{ D = 1; E = 1; H = 0 }
}

$colIndex = @{ D = 4; E = 5; F = 6; G = 7; H = 8 }

foreach ($r in $values.Keys) {
    $rowVals = $values[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Cells.Item([int]$r, $colIndex[$col]).Value = $rowVals[$col]
    }
}
